$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the header text in E1 ("Ingreso 1") so the column is removed from the dataset
$ws.Range("E1").ClearContents()

# Move the active selection to J7, matching the cell the user left selected when saving
$ws.Range("J7").Select()
